$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (0..18), columns A..S
# Row 1
$ws.Cells.Item(1, 1).Value = 0
$ws.Cells.Item(1, 2).Value = 1
$ws.Cells.Item(1, 3).Value = 2
$ws.Cells.Item(1, 4).Value = 3
$ws.Cells.Item(1, 5).Value = 4
$ws.Cells.Item(1, 6).Value = 5
$ws.Cells.Item(1, 7).Value = 6
$ws.Cells.Item(1, 8).Value = 7
$ws.Cells.Item(1, 9).Value = 8
$ws.Cells.Item(1, 10).Value = 9
$ws.Cells.Item(1, 11).Value = 10
$ws.Cells.Item(1, 12).Value = 11
$ws.Cells.Item(1, 13).Value = 12
$ws.Cells.Item(1, 14).Value = 13
$ws.Cells.Item(1, 15).Value = 14
$ws.Cells.Item(1, 16).Value = 15
$ws.Cells.Item(1, 17).Value = 16
$ws.Cells.Item(1, 18).Value = 17
$ws.Cells.Item(1, 19).Value = 18

# Row 2
$ws.Cells.Item(2, 1).Value = -1.542905112095295
$ws.Cells.Item(2, 2).Value = -107.6723250834219
$ws.Cells.Item(2, 3).Value = 117.5942845603324
$ws.Cells.Item(2, 4).Value = -2078216.403033835
$ws.Cells.Item(2, 5).Value = -722530669.6681616
$ws.Cells.Item(2, 6).Value = -218835220559.0445
$ws.Cells.Item(2, 7).Value = -63487359319802.37
$ws.Cells.Item(2, 8).Value = [double]"-1.7888376881237e+16"
$ws.Cells.Item(2, 9).Value = [double]"-4.910387891247182e+18"
$ws.Cells.Item(2, 10).Value = [double]"-1.314956972823998e+21"
$ws.Cells.Item(2, 11).Value = [double]"-3.435029366888184e+23"
$ws.Cells.Item(2, 12).Value = [double]"-8.769151142633699e+25"
$ws.Cells.Item(2, 13).Value = [double]"-2.199861395007915e+28"
$ws.Cells.Item(2, 14).Value = [double]"-5.464472141685351e+30"
$ws.Cells.Item(2, 15).Value = [double]"-1.349188375047402e+33"
$ws.Cells.Item(2, 16).Value = [double]"-3.315948831958216e+35"
$ws.Cells.Item(2, 17).Value = [double]"-8.12761329178604e+37"
$ws.Cells.Item(2, 18).Value = [double]"-1.990483672455331e+40"
$ws.Cells.Item(2, 19).Value = [double]"-4.875869488341349e+42"

# Row 3
$ws.Cells.Item(3, 1).Value = -0.5187822358978966
$ws.Cells.Item(3, 2).Value = -96.26250275744302
$ws.Cells.Item(3, 3).Value = -21774.18077738656
$ws.Cells.Item(3, 4).Value = -3939452.718682132
$ws.Cells.Item(3, 5).Value = -743370511.1757892
$ws.Cells.Item(3, 6).Value = -145173280549.7003
$ws.Cells.Item(3, 7).Value = -28532797857823.7
$ws.Cells.Item(3, 8).Value = -5381414459277003
$ws.Cells.Item(3, 9).Value = [double]"1.005384841232984e+18"
$ws.Cells.Item(3, 10).Value = [double]"1.846293576704401e+20"
$ws.Cells.Item(3, 11).Value = [double]"3.296553628287555e+22"
$ws.Cells.Item(3, 12).Value = [double]"5.623688490542846e+24"
$ws.Cells.Item(3, 13).Value = [double]"8.891934667463711e+26"
$ws.Cells.Item(3, 14).Value = [double]"1.232374983626391e+29"
$ws.Cells.Item(3, 15).Value = [double]"1.331724709637454e+31"
$ws.Cells.Item(3, 16).Value = [double]"7.180190329923423e+32"
$ws.Cells.Item(3, 17).Value = [double]"-9.786830923925587e+34"
$ws.Cells.Item(3, 18).Value = [double]"-3.417179406319799e+37"
$ws.Cells.Item(3, 19).Value = [double]"-6.065716220636713e+39"

# Row 4
$ws.Cells.Item(4, 1).Value = -0.9715511871506971
$ws.Cells.Item(4, 2).Value = -47.91540311645265
$ws.Cells.Item(4, 3).Value = 5333.121447522998
$ws.Cells.Item(4, 4).Value = 899618.7245339525
$ws.Cells.Item(4, 5).Value = 142790863.9063835
$ws.Cells.Item(4, 6).Value = 19892354175.39549
$ws.Cells.Item(4, 7).Value = 2193980242931.121
$ws.Cells.Item(4, 8).Value = 356656788858720.1
$ws.Cells.Item(4, 9).Value = [double]"6.159645247031888e+16"
$ws.Cells.Item(4, 10).Value = [double]"1.170956605483629e+19"
$ws.Cells.Item(4, 11).Value = [double]"2.357134820742203e+21"
$ws.Cells.Item(4, 12).Value = [double]"4.628582281339793e+23"
$ws.Cells.Item(4, 13).Value = [double]"8.10522661824741e+25"
$ws.Cells.Item(4, 14).Value = [double]"1.162721903641395e+28"
$ws.Cells.Item(4, 15).Value = [double]"1.360689617654482e+30"
$ws.Cells.Item(4, 16).Value = [double]"2.088608903231858e+32"
$ws.Cells.Item(4, 17).Value = [double]"3.801299579479052e+34"
$ws.Cells.Item(4, 18).Value = [double]"7.316263348866375e+36"
$ws.Cells.Item(4, 19).Value = [double]"1.455426636667734e+39"

# Row 5
$ws.Cells.Item(5, 1).Value = -0.7138522768094833
$ws.Cells.Item(5, 2).Value = 79.20981309627001
$ws.Cells.Item(5, 3).Value = 9275.66480615844
$ws.Cells.Item(5, 4).Value = 1396088.291332924
$ws.Cells.Item(5, 5).Value = 276334716.5663702
$ws.Cells.Item(5, 6).Value = 56783985270.84228
$ws.Cells.Item(5, 7).Value = 11385693824000.71
$ws.Cells.Item(5, 8).Value = 2101600583587585
$ws.Cells.Item(5, 9).Value = [double]"3.650676638760862e+17"
$ws.Cells.Item(5, 10).Value = [double]"5.770318972107194e+19"
$ws.Cells.Item(5, 11).Value = [double]"-7.862388380938082e+21"
$ws.Cells.Item(5, 12).Value = [double]"-8.607003515235669e+23"
$ws.Cells.Item(5, 13).Value = [double]"-7.427121555094514e+25"
$ws.Cells.Item(5, 14).Value = [double]"-7.955168817271931e+27"
$ws.Cells.Item(5, 15).Value = [double]"-1.377268362418268e+30"
$ws.Cells.Item(5, 16).Value = [double]"-2.386224608982165e+32"
$ws.Cells.Item(5, 17).Value = [double]"-4.22328746981304e+34"
$ws.Cells.Item(5, 18).Value = [double]"-7.703688064493114e+36"
$ws.Cells.Item(5, 19).Value = [double]"-1.42965178751053e+39"

# Row 6
$ws.Cells.Item(6, 1).Value = 0.8231336624746795
$ws.Cells.Item(6, 2).Value = 36.79575802539394
$ws.Cells.Item(6, 3).Value = 341.7735292848716
$ws.Cells.Item(6, 4).Value = -566942.3603989104
$ws.Cells.Item(6, 5).Value = -64507535.7913077
$ws.Cells.Item(6, 6).Value = -6368476291.501096
$ws.Cells.Item(6, 7).Value = -476457990844.7783
$ws.Cells.Item(6, 8).Value = 33050504344210.38
$ws.Cells.Item(6, 9).Value = [double]"1.380952373738552e+16"
$ws.Cells.Item(6, 10).Value = [double]"1.990339438062695e+18"
$ws.Cells.Item(6, 11).Value = [double]"2.404301331440669e+20"
$ws.Cells.Item(6, 12).Value = [double]"2.970316171607975e+22"
$ws.Cells.Item(6, 13).Value = [double]"3.45682402395027e+24"
$ws.Cells.Item(6, 14).Value = [double]"4.173607128946785e+26"
$ws.Cells.Item(6, 15).Value = [double]"5.13039575239473e+28"
$ws.Cells.Item(6, 16).Value = [double]"6.213334344196775e+30"
$ws.Cells.Item(6, 17).Value = [double]"8.009580346051177e+32"
$ws.Cells.Item(6, 18).Value = [double]"1.20502165607286e+35"
$ws.Cells.Item(6, 19).Value = [double]"2.029777869283875e+37"

# Apply header style (s="1") to the newly added header cells C1:S1
$ws.Range("B1").Copy()
$ws.Range("C1:S1").PasteSpecial(-4122)
$excel.CutCopyMode = 0